$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.517.72'
$ws.Range("E2").Value = '  -0.40%  '

$ws.Range("D3").Value = '2.268.76'
$ws.Range("E3").Value = '  -0.58%  '

$ws.Range("E4").Value = '  -0.17%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '119.02'
$ws.Range("E5").Value = '  +4.37%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '265.59'
$ws.Range("E6").Value = '  -0.72%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.645'
$ws.Range("E7").Value = '  +3.23%  '

$ws.Range("E8").Value = '  +0.24%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.620'
$ws.Range("E9").Value = '  +1.43%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '47.64'
$ws.Range("E10").Value = '  -1.95%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0945'
$ws.Range("E11").Value = '  +0.72%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '9.18'
$ws.Range("E12").Value = '  +3.92%  '

$ws.Range("E13").Value = '  -1.39%  '

$ws.Range("E14").Value = '  -2.12%  '

$ws.Range("E15").Value = '  +3.10%  '

$ws.Range("D16").Value = '2.611.11'
$ws.Range("E16").Value = '  -0.47%  '

$ws.Range("D17").Value = '2.265.65'
$ws.Range("E17").Value = '  -0.58%  '

$ws.Range("D18").Value = '43.538.75'
$ws.Range("E18").Value = '  +0.09%  '

$ws.Range("E19").Value = '  +0.73%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.84'
$ws.Range("E20").Value = '  -3.11%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.08'
$ws.Range("E21").Value = '  +0.10%  '

$ws.Range("E22").Value = '  -0.49%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.55'
$ws.Range("E23").Value = '  +1.29%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.50'
$ws.Range("E24").Value = '  -5.05%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.87'
$ws.Range("E25").Value = '  -0.20%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.01'
$ws.Range("E26").Value = '  +4.11%  '

$ws.Range("E27").Value = '  +1.80%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '41.60'
$ws.Range("E28").Value = '  +1.27%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.39'
$ws.Range("E29").Value = '  -0.32%  '

$ws.Range("E30").Value = '  -0.27%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '172.10'
$ws.Range("E31").Value = '  -0.82%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.57'
$ws.Range("E32").Value = '  +0.34%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0913'
$ws.Range("E33").Value = '  +0.32%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.70'
$ws.Range("E34").Value = '  +0.83%  '

$ws.Range("E35").Value = '  +2.40%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0384'
$ws.Range("E36").Value = '  +9.45%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.18'
$ws.Range("E37").Value = '  +11.39%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.55'
$ws.Range("E38").Value = '  -1.42%  '

$ws.Range("E39").Value = '  +0.95%  '

$ws.Range("E40").Value = '  +5.07%  '

$ws.Range("B41").Value = 'MultiversX'
$ws.Range("C41").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '73.78'
$ws.Range("E41").Value = '  -1.61%  '

$ws.Range("B42").Value = 'Celestia'
$ws.Range("C42").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '13.74'
$ws.Range("E42").Value = '  -4.32%  '

$ws.Range("E43").Value = '  -0.57%  '

$ws.Range("E44").Value = '  -0.12%  '

$ws.Range("E45").Value = '  -1.03%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.79'
$ws.Range("E46").Value = '  -6.62%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '73.69'
$ws.Range("E47").Value = '  +40.44%  '

$ws.Range("E48").Value = '  +1.51%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.54'

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0998'
$ws.Range("E50").Value = '  +0.29%  '
